$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44214

$ws.Range("D3").Value = 44855
$ws.Range("J3").Value = 7900

$ws.Range("D4").Value = 44602
$ws.Range("J4").Value = 12000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = 'Provincia de Chacabuco'
$ws.Range("P4").Value = 30

$ws.Range("D5").Value = 44602
$ws.Range("I5").Value = 'Segunda'
$ws.Range("J5").Value = 6000
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 25

$ws.Range("D6").Value = 44167
$ws.Range("J6").Value = 7000

$ws.Range("D7").Value = 44875

$ws.Range("D8").Value = 44186
$ws.Range("J8").Value = 10000

$ws.Range("D10").Value = 44847
$ws.Range("J10").Value = 7900

$ws.Range("D11").Value = 44231
$ws.Range("J11").Value = 12000

$ws.Range("D12").Value = 44162

$ws.Range("D13").Value = 44860

$ws.Range("D14").Value = 44845
$ws.Range("J14").Value = 7900

$ws.Range("D15").Value = 44901
$ws.Range("J15").Value = 7000

$ws.Range("D16").Value = 44881
$ws.Range("J16").Value = 7900

$ws.Range("D17").Value = 44181

$ws.Range("D18").Value = 44189
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 16000
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("P18").Value = 30

$ws.Range("D19").Value = 44166
$ws.Range("J19").Value = 7000

$ws.Range("D20").Value = 44209
$ws.Range("J20").Value = 7000

$ws.Range("D21").Value = 44168
$ws.Range("J21").Value = 7000

$ws.Range("D22").Value = 44902
$ws.Range("J22").Value = 7000

$ws.Range("D23").Value = 44600
$ws.Range("J23").Value = 1300
$ws.Range("K23").Value = 3500
$ws.Range("L23").Value = 4000
$ws.Range("M23").Value = 3808
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 38

$ws.Range("D24").Value = 44204
$ws.Range("J24").Value = 7000

$ws.Range("D25").Value = 44215
$ws.Range("J25").Value = 16000

$ws.Range("D26").Value = 44161
$ws.Range("J26").Value = 7000

$ws.Range("D27").Value = 44188
$ws.Range("J27").Value = 12000

$ws.Range("D28").Value = 44882
$ws.Range("J28").Value = 7900
$ws.Range("K28").Value = 3000
$ws.Range("M28").Value = 3000
$ws.Range("P28").Value = 30

$ws.Range("D29").Value = 44210
$ws.Range("J29").Value = 8800
$ws.Range("K29").Value = 2500
$ws.Range("M29").Value = 2750
$ws.Range("O29").Value = 'Provincia de Chacabuco'
$ws.Range("P29").Value = 28

$ws.Range("D30").Value = 44187
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 12000
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = 3000
$ws.Range("O30").Value = 'Provincia de Chacabuco'
$ws.Range("P30").Value = 30

$ws.Range("D31").Value = 44846
$ws.Range("J31").Value = 7900

$ws.Range("D32").Value = 44229
$ws.Range("J32").Value = 16000

$ws.Range("D33").Value = 44873
$ws.Range("J33").Value = 7900

$ws.Range("D34").Value = 44160

$ws.Range("D35").Value = 44230
$ws.Range("J35").Value = 16000

$ws.Range("D36").Value = 44880
$ws.Range("J36").Value = 7900

$ws.Range("D37").Value = 44159

$ws.Range("D38").Value = 44874

$ws.Range("D39").Value = 44245
$ws.Range("J39").Value = 9000
$ws.Range("O39").Value = 'Región Metropolitana'

$ws.Range("D40").Value = 44245
$ws.Range("I40").Value = 'Segunda'
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = 2500
$ws.Range("O40").Value = 'Región Metropolitana'
$ws.Range("P40").Value = 25

$ws.Range("D41").Value = 44859
$ws.Range("J41").Value = 7900

$ws.Range("D42").Value = 44883
$ws.Range("J42").Value = 9700
